$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values in column D (e.g. "297.96", "6.00", "0.770") look like
# plain decimal numbers. Excel auto-converts such text to a real number when
# assigned through COM, which would silently normalize away significant
# trailing zeros/formatting ("6.00" -> 6, "0.770" -> 0.77, "167.10" -> 167.1,
# ...). Marking those specific cells as Text first keeps the new price string
# exactly as published. Prices that contain more than one "." (e.g.
# "42.057.77") can never be parsed as a number anyway, and the "Volume(1h)"
# column E values (padded with spaces and a trailing "%") are never number-like
# either, so neither needs this treatment.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.057.77"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").Value = "2.262.48"
$ws.Range("E3").Value = "  -3.08%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "297.96"
$ws.Range("E5").Value = "  -2.68%  "

$ws.Range("D6").Value = "93.69"
$ws.Range("E6").Value = "  -7.50%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -4.05%  "

$ws.Range("D10").Value = "32.89"
$ws.Range("E10").Value = "  -6.06%  "

$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("D12").Value = "48.13"
$ws.Range("E12").Value = "  -7.79%  "

$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").Value = "6.63"
$ws.Range("E14").Value = "  -2.65%  "

$ws.Range("D15").Value = "2.614.62"
$ws.Range("E15").Value = "  -3.28%  "

$ws.Range("D16").Value = "15.49"
$ws.Range("E16").Value = "  -2.35%  "

$ws.Range("D17").Value = "2.267.56"
$ws.Range("E17").Value = "  -4.25%  "

$ws.Range("D18").Value = "0.770"
$ws.Range("E18").Value = "  -5.02%  "

$ws.Range("D19").Value = "42.056.46"
$ws.Range("E19").Value = "  -1.89%  "

$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("E20").Value = "  -2.41%  "

$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  -3.23%  "

$ws.Range("D22").Value = "11.32"
$ws.Range("E22").Value = "  -3.14%  "

$ws.Range("D23").Value = "66.56"
$ws.Range("E23").Value = "  -1.96%  "

$ws.Range("D24").Value = "232.76"
$ws.Range("E24").Value = "  -1.65%  "

$ws.Range("E25").Value = "  -4.56%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -4.41%  "

$ws.Range("D28").Value = "23.73"
$ws.Range("E28").Value = "  -6.53%  "

$ws.Range("E29").Value = "  -1.03%  "

$ws.Range("D30").Value = "167.10"
$ws.Range("E30").Value = "  +3.79%  "

$ws.Range("D31").Value = "33.50"
$ws.Range("E31").Value = "  -4.11%  "

$ws.Range("D32").Value = "9.01"
$ws.Range("E32").Value = "  -3.93%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "4.90"
$ws.Range("E34").Value = "  -4.23%  "

$ws.Range("E35").Value = "  -3.26%  "

$ws.Range("E36").Value = "  -5.65%  "

$ws.Range("D37").Value = "0.0688"
$ws.Range("E37").Value = "  -5.37%  "

$ws.Range("D38").Value = "16.10"
$ws.Range("E38").Value = "  -7.88%  "

$ws.Range("D39").Value = "2.76"
$ws.Range("E39").Value = "  -5.69%  "

$ws.Range("E40").Value = "  -3.76%  "

$ws.Range("E41").Value = "  -3.80%  "

$ws.Range("E42").Value = "  -8.50%  "

$ws.Range("E43").Value = "  -1.87%  "

$ws.Range("D44").Value = "1.954.95"
$ws.Range("E44").Value = "  -2.46%  "

$ws.Range("E45").Value = "  -2.76%  "

$ws.Range("D46").Value = "17.17"
$ws.Range("E46").Value = "  -8.40%  "

$ws.Range("E47").Value = "  -6.51%  "

$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  -5.53%  "

$ws.Range("E49").Value = "  -3.16%  "

$ws.Range("D50").Value = "2.488.01"
$ws.Range("E50").Value = "  -2.83%  "

$ws.Range("D51").Value = "51.53"
$ws.Range("E51").Value = "  -7.37%  "
